$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    existing "2021-Q1" sheet), so the tab order becomes:
#    总计, 2022-Q3, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# NOTE: worksheet references returned by `.Item(index)` track the *slot*,
# not the sheet identity, so fetch "2021-Q1" by name now that the new sheet
# has shifted everyone after "总计" down by one slot.
$q1sheet = $wb.Worksheets.Item("2021-Q1")

# Pull header formatting (bold / border / alignment style) from the
# "2021-Q1" sheet so the new sheet's header row matches the house style.
$q1sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q1sheet.Range("A2:A3").Copy($newSheet.Range("A2:A3"))

# Header row text (column D differs: "基金规模" instead of "基金金额")
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Numeric-looking fund codes / figures must stay text, like the sibling
# sheets (013009 would otherwise lose its leading zero).
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "013009"
$newSheet.Range("C2").Value = "万家港股通精选混合A"
$newSheet.Range("D2").Value = "1.78"
$newSheet.Range("E2").Value = "87.54"
$newSheet.Range("F2").Value = "3.76"
$newSheet.Range("G2").Value = "0.0669"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "013010"
$newSheet.Range("C3").Value = "万家港股通精选混合C"
$newSheet.Range("D3").Value = "0.71"
$newSheet.Range("E3").Value = "87.54"
$newSheet.Range("F3").Value = "3.76"
$newSheet.Range("G3").Value = "0.0267"
$newSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new row for 2022-Q3 as row 2,
#    pushing 2021-Q1 -> row 3 and 2020-Q4 -> row 4, and renumber the index
#    column (A) to 0,1,2.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.09

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q1"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 1.65

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2020-Q4"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 1.04

# ---------------------------------------------------------------------------
# 3) Adding the new sheet shifts the active tab onto it; restore the
#    original selection ("2020-Q4", the last sheet) so that sheet keeps its
#    tabSelected state, unchanged by this edit.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()

